$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the two new columns (I: I0, J: IF)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the look of the existing header cells (bold, centered, bordered)
$ws.Range("I1:J1").Font.Bold = $true
$ws.Range("I1:J1").HorizontalAlignment = -4108
$ws.Range("I1:J1").VerticalAlignment = -4160
$ws.Range("I1:J1").Borders.LineStyle = 1

# New data for columns I (I0) and J (IF), rows 2-48
$ijValues = @{
    "2" = @(7, 8)
    "3" = @(7, 8)
    "4" = @(8, 9)
    "5" = @(8, 8)
    "6" = @(7, 9)
    "7" = @(5, 6)
    "8" = @(7, 8)
    "9" = @(3, 6)
    "10" = @(3, 5)
    "11" = @(1, 2)
    "12" = @(1, 1)
    "13" = @(5, 6)
    "14" = @(5, 6)
    "15" = @(6, 7)
    "16" = @(6, 7)
    "17" = @(5, 6)
    "18" = @(9, 10)
    "19" = @(7, 7)
    "20" = @(6, 6)
    "21" = @(9, 9)
    "22" = @(7, 7)
    "23" = @(1, 1)
    "24" = @(7, 7)
    "25" = @(4, 4)
    "26" = @(8, 8)
    "27" = @(6, 6)
    "28" = @(7, 7)
    "29" = @(9, 9)
    "30" = @(6, 6)
    "31" = @(5, 6)
    "32" = @(7, 7)
    "33" = @(7, 8)
    "34" = @(8, 9)
    "35" = @(8, 8)
    "36" = @(10, 10)
    "37" = @(5, 7)
    "38" = @(7, 8)
    "39" = @(4, 6)
    "40" = @(6, 7)
    "41" = @(7, 7)
    "42" = @(5, 6)
    "43" = @(8, 8)
    "44" = @(8, 8)
    "45" = @(3, 4)
    "46" = @(5, 5)
    "47" = @(6, 6)
    "48" = @(6, 6)
}

foreach ($row in $ijValues.Keys) {
    $pair = $ijValues[$row]
    $ws.Cells.Item([int]$row, 9).Value = $pair[0]
    $ws.Cells.Item([int]$row, 10).Value = $pair[1]
}

